$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.450.68'
$ws.Range("E2").Value = '  +4.22%  '
$ws.Range("D3").Value = '1.804.29'
$ws.Range("E3").Value = '  +1.69%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.19'
$ws.Range("E5").Value = '  +0.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5503'
$ws.Range("E7").Value = '  +4.74%  '
$ws.Range("E8").Value = '  +6.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07615'
$ws.Range("E9").Value = '  +3.45%  '
$ws.Range("E10").Value = '  -0.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.128'
$ws.Range("E11").Value = '  +3.41%  '
$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.24'
$ws.Range("E13").Value = '  +3.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.199'
$ws.Range("E14").Value = '  +2.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.475'
$ws.Range("E15").Value = '  +7.44%  '
$ws.Range("D16").Value = '1.807.60'
$ws.Range("E16").Value = '  +2.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.09'
$ws.Range("E17").Value = '  +3.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001073'
$ws.Range("E18").Value = '  +2.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06443'
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9996'
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.35'
$ws.Range("E21").Value = '  +3.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.977'
$ws.Range("E22").Value = '  +2.68%  '
$ws.Range("D23").Value = '28.448.37'
$ws.Range("E23").Value = '  +3.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.43'
$ws.Range("E24").Value = '  +1.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.140'
$ws.Range("E25").Value = '  +2.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.04'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.72'
$ws.Range("E27").Value = '  +2.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.415'
$ws.Range("E28").Value = '  +3.25%  '
$ws.Range("D29").Value = '2.012.01'
$ws.Range("E29").Value = '  +1.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.95'
$ws.Range("E30").Value = '  +2.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.126'
$ws.Range("E31").Value = '  +6.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1022'
$ws.Range("E32").Value = '  +4.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.788'
$ws.Range("E33").Value = '  +4.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.688'
$ws.Range("E34").Value = '  +1.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2318'
$ws.Range("E35").Value = '  +14.73%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06459'
$ws.Range("E36").Value = '  +8.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02328'
$ws.Range("E37").Value = '  +4.64%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.181'
$ws.Range("E38").Value = '  +7.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.832'
$ws.Range("E39").Value = '  +9.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.68'
$ws.Range("E40").Value = '  +4.35%  '
$ws.Range("E41").Value = '  +4.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.164'
$ws.Range("E42").Value = '  +1.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9998'
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.385'
$ws.Range("E44").Value = '  -3.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.62'
$ws.Range("E45").Value = '  +3.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5991'
$ws.Range("E46").Value = '  +4.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.681'
$ws.Range("E47").Value = '  +1.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.21'
$ws.Range("E48").Value = '  +5.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.985'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.152'
$ws.Range("E50").Value = '  +3.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06893'
$ws.Range("E51").Value = '  +2.72%  '
